# Diplomarbeit Arbeitszeit Pichler - "download faulty, working on fix"
# Append new time-tracking entries (rows 41-44) to the log table on Tabelle1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New rows appended at the bottom of the Datum/Zeit/Einheit/Tätigkeit/Projekt table (cols E:I)
# Use the same raw date serials (whole-day, no time-of-day component) as the
# rest of the E column, and copy the existing date cell's style so no new
# number-format style gets created.
$ws.Range("E41").Value = 43697
$ws.Range("F41").Value = 4
$ws.Range("G41").Value = "Stunden"
$ws.Range("H41").Value = "Programmieren"
$ws.Range("I41").Value = "Theorie"

$ws.Range("E42").Value = 43698
$ws.Range("F42").Value = 3
$ws.Range("G42").Value = "Stunden"
$ws.Range("H42").Value = "Programmieren"
$ws.Range("I42").Value = "Logging optimiert, Teil 1"

$ws.Range("E43").Value = 43700
$ws.Range("F43").Value = 4
$ws.Range("G43").Value = "Stunden"
$ws.Range("H43").Value = "Programmieren"

$ws.Range("E44").Value = 43701
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = "Stunden"
$ws.Range("H44").Value = "Nichts"

# Match date formatting used by the rest of the column (style of E40)
$ws.Range("E40").Copy() | Out-Null
$ws.Range("E41:E44").PasteSpecial(-4122) | Out-Null

# Update the active selection to reflect the newly-added rows at the
# bottom of the sheet.
$ws.Range("H44").Select()
